$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 911.64514
$ws.Range("J17").Value = 925.4
$ws.Range("L17").Value = 2776.2
$ws.Range("N17").Value = -3112.2
$ws.Range("H61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H76").Value = 3639.75
$ws.Range("I76").Value = 3588.2856
$ws.Range("K76").Value = 3588.2856
$ws.Range("M76").Value = -3273.2856
$ws.Range("H79").Value = 3639.75
$ws.Range("I79").Value = 3588.2856
$ws.Range("K79").Value = 3588.2856
$ws.Range("M79").Value = -2496.2856
$ws.Range("H100").Value = 2021.4
$ws.Range("I100").Value = 2051
$ws.Range("J100").Value = 1903
$ws.Range("K100").Value = 2051
$ws.Range("L100").Value = 1903
$ws.Range("M100").Value = -1510
$ws.Range("N100").Value = -2985
$ws.Range("H111").Value = 1300
$ws.Range("I111").Value = 1300
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3900
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -833
$ws.Range("H113").Value = 30306832
$ws.Range("J113").Value = 5053.6
$ws.Range("L113").Value = 5053.6
$ws.Range("N113").Value = -11561.6
$ws.Range("H131").Value = 1012836
$ws.Range("I131").Value = 1444994.6
$ws.Range("K131").Value = 4334983.800000001
$ws.Range("M131").Value = -4329943.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1132.6923
$ws.Range("I2").Value = 969.8889
$ws.Range("J2").Value = 1499
$ws.Range("K2").Value = 969.8889
$ws.Range("L2").Value = 1499
$ws.Range("M2").Value = -856.8889
$ws.Range("N2").Value = -1725
$ws.Range("H45").Value = 24310.947
$ws.Range("I45").Value = 40072.816
$ws.Range("K45").Value = 40072.816
$ws.Range("M45").Value = -39695.816
$ws.Range("H61").Value = 2985.1072
$ws.Range("I61").Value = 2777.5293
$ws.Range("K61").Value = 2777.5293
$ws.Range("M61").Value = -2565.5293
$ws.Range("H97").Value = 1020.7273
$ws.Range("I97").Value = 1009.8095
$ws.Range("K97").Value = 1009.8095
$ws.Range("M97").Value = -513.8095
$ws.Range("H116").Value = 1132.6923
$ws.Range("I116").Value = 969.8889
$ws.Range("J116").Value = 1499
$ws.Range("K116").Value = 969.8889
$ws.Range("L116").Value = 1499
$ws.Range("M116").Value = 1324.1111
$ws.Range("N116").Value = -6087
$ws.Range("H136").Value = 2985.1072
$ws.Range("I136").Value = 2777.5293
$ws.Range("K136").Value = 8332.5879
$ws.Range("M136").Value = -5782.5879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1132.6923
$ws.Range("I3").Value = 969.8889
$ws.Range("J3").Value = 1499
$ws.Range("K3").Value = 969.8889
$ws.Range("L3").Value = 1499
$ws.Range("M3").Value = -855.8889
$ws.Range("N3").Value = -1727
$ws.Range("H86").Value = 3205.2144
$ws.Range("J86").Value = 4375
$ws.Range("L86").Value = 4375
$ws.Range("N86").Value = -6621
$ws.Range("H89").Value = 3205.2144
$ws.Range("J89").Value = 4375
$ws.Range("L89").Value = 21875
$ws.Range("N89").Value = -33107
$ws.Range("H103").Value = 30399.8
$ws.Range("J103").Value = 30399.8
$ws.Range("L103").Value = 30399.8
$ws.Range("N103").Value = -32743.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4036990.2
$ws.Range("I31").Value = 2755.2727
$ws.Range("J31").Value = 6255819.5
$ws.Range("K31").Value = 2755.2727
$ws.Range("L31").Value = 6255819.5
$ws.Range("M31").Value = -2460.2727
$ws.Range("N31").Value = -6256409.5
$ws.Range("H34").Value = 4036990.2
$ws.Range("I34").Value = 2755.2727
$ws.Range("J34").Value = 6255819.5
$ws.Range("K34").Value = 2755.2727
$ws.Range("L34").Value = 6255819.5
$ws.Range("M34").Value = -2553.2727
$ws.Range("N34").Value = -6256223.5
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H132").Value = 10421525
$ws.Range("I132").Value = 3537.1667
$ws.Range("K132").Value = 10611.5001
$ws.Range("M132").Value = -8081.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 700
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 132
$ws.Range("I61").Value = 132
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 396
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -181
$ws.Range("H98").Value = 618
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 476.30304
$ws.Range("J107").Value = 466.2069
$ws.Range("L107").Value = 1398.6207
$ws.Range("N107").Value = -5238.620699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 22547222
$ws.Range("I11").Value = 56875000
$ws.Range("J11").Value = 5383333.5
$ws.Range("K11").Value = 56875000
$ws.Range("L11").Value = 5383333.5
$ws.Range("M11").Value = -56874861
$ws.Range("N11").Value = -5383611.5
$ws.Range("H132").Value = 2533.7437
$ws.Range("I132").Value = 1935.9474
$ws.Range("J132").Value = 3101.65
$ws.Range("K132").Value = 5807.8422
$ws.Range("L132").Value = 9304.950000000001
$ws.Range("M132").Value = -3277.8422
$ws.Range("N132").Value = -14364.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1344.6666
$ws.Range("I10").Value = 360.6
$ws.Range("J10").Value = 2574.75
$ws.Range("K10").Value = 360.6
$ws.Range("L10").Value = 2574.75
$ws.Range("M10").Value = -220.6
$ws.Range("N10").Value = -2854.75
$ws.Range("H13").Value = 4002.75
$ws.Range("I13").Value = 3005.5
$ws.Range("K13").Value = 3005.5
$ws.Range("M13").Value = -2865.5
$ws.Range("H46").Value = 2825.2307
$ws.Range("I46").Value = 1975.3636
$ws.Range("J46").Value = 7499.5
$ws.Range("K46").Value = 1975.3636
$ws.Range("L46").Value = 7499.5
$ws.Range("M46").Value = -1787.3636
$ws.Range("N46").Value = -7875.5
$ws.Range("H104").Value = 35472.4
$ws.Range("J104").Value = 35472.4
$ws.Range("L104").Value = 35472.4
$ws.Range("N104").Value = -42460.4
$ws.Range("H132").Value = 4375.9614
$ws.Range("I132").Value = 4236.8125
$ws.Range("J132").Value = 4598.6
$ws.Range("K132").Value = 12710.4375
$ws.Range("L132").Value = 13795.8
$ws.Range("M132").Value = -10180.4375
$ws.Range("N132").Value = -18855.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 51857.25
$ws.Range("J46").Value = 51857.25
$ws.Range("L46").Value = 51857.25
$ws.Range("N46").Value = -52319.25
$ws.Range("H101").Value = 39911
$ws.Range("J101").Value = 39911
$ws.Range("L101").Value = 39911
$ws.Range("N101").Value = -46401
$ws.Range("H110").Value = 115000
$ws.Range("J110").Value = 115000
$ws.Range("L110").Value = 115000
$ws.Range("N110").Value = -123180
$ws.Range("H134").Value = 51857.25
$ws.Range("J134").Value = 51857.25
$ws.Range("L134").Value = 155571.75
$ws.Range("N134").Value = -160641.75
